# [ CHECK FRAME ID BEFORE ADD SIGNAL ATTRIBUTES ]
# On the K-Matrix sheet, column K holds "Message Type". For every data row
# the value is updated from "/" to "STD".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("K-Matrix ")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 13 }

for ($row = 2; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 11).Value = "STD"
}
